$wb = $excel.ActiveWorkbook
$userSheet = $wb.Worksheets.Item("User")

# Add the "Evaluation Warning" sheet right after the "User" sheet and make
# it the active sheet, mirroring the Aspose.Cells unlicensed-mode banner.
$warnSheet = $wb.Worksheets.Add([System.Type]::Missing, $userSheet)
$warnSheet.Name = "Evaluation Warning"

$cell = $warnSheet.Range("A5")
$cell.Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2020 Aspose Pty Ltd."
$cell.Font.Name = "Arial"
$cell.Font.Size = 18
$cell.Font.Bold = $true
$cell.Font.Italic = $true
$cell.Font.Color = 16711680

$warnSheet.Rows.Item(5).RowHeight = 23.25
